# Updating version to 2 — see changelog.
#
# "results" sheet: drop the second sample run (row 3), widen the header
# row with the new "-UN" variant columns (N:S), and refresh the row-2
# measurements so every solver column (H:J) now reports a real elapsed
# time instead of a boolean placeholder.
#
# "stats" sheet: drop the middle "run 1" block entirely, and extend each
# remaining block (run 0 / Average) with the three new "-UN" variant rows
# before the trailing Kruskal baseline row. Re-merge the A column labels
# to span the new, longer blocks.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet "results"
# ---------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("results")

# Drop the "run 1" sample row (old row 3) first, so row 2 is the only
# data row left.
$ws1.Rows.Item(3).Delete()

# The single "-unmerged" column (G) is replaced by 4 new "-UN" variant
# columns, which pushes every following header right by 3 slots; then
# the "-unmergedND" column is likewise replaced by 4 "-UN...ND" columns
# in the newly-appended N:S range. Re-write the full header label set
# (G1:S1) so the text matches the new, longer column layout.
$ws1.Range("G1").Value = "S*-BS-UN"
$ws1.Range("H1").Value = "S*-HS-UN"
$ws1.Range("I1").Value = "S*-MM-UN"
$ws1.Range("J1").Value = "S*-MM0-UN"
$ws1.Range("K1").Value = "var"
$ws1.Range("L1").Value = "S*-BSND"
$ws1.Range("M1").Value = "S*-HSND"
$ws1.Range("N1").Value = "S*-MMND"
$ws1.Range("O1").Value = "S*-MM0ND"
$ws1.Range("P1").Value = "S*-BS-UNND"
$ws1.Range("Q1").Value = "S*-HS-UNND"
$ws1.Range("R1").Value = "S*-MM-UNND"
$ws1.Range("S1").Value = "S*-MM0-UNND"

# Apply the same header style (bold, centered, boxed) used by the rest
# of row 1 to the newly-added N:S cells.
$ws1.Range("N1:S1").Style = $ws1.Range("M1").Style

# Refresh the measurements on row 2 with the new benchmark numbers.
$ws1.Range("B2").Value = 270.496
$ws1.Range("C2").Value = 270.496
$ws1.Range("D2").Value = 270.496
$ws1.Range("E2").Value = 270.496
$ws1.Range("F2").Value = 270.496
$ws1.Range("G2").Value = 270.496
$ws1.Range("H2").Value = 270.496
$ws1.Range("I2").Value = 270.496
$ws1.Range("J2").Value = 270.496
$ws1.Range("K2").Value = 0

# New boolean placeholder cells for the "-UN" variant columns.
$ws1.Range("N2").Value = $false
$ws1.Range("O2").Value = $false
$ws1.Range("P2").Value = $false
$ws1.Range("Q2").Value = $false
$ws1.Range("R2").Value = $false
$ws1.Range("S2").Value = $false

# ---------------------------------------------------------------------
# Sheet "stats"
# ---------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("stats")

# Unmerge the three existing per-run label blocks before restructuring
# the rows underneath them.
$ws2.Range("A2:A7").UnMerge()
$ws2.Range("A8:A13").UnMerge()
$ws2.Range("A14:A19").UnMerge()

# Drop the middle "run 1" block (old rows 8-13) entirely.
$ws2.Range("A8:K13").Delete()

# Insert 3 fresh rows before the trailing "Kruskal" row in each of the
# two remaining blocks, to hold the new "-UN" variant entries.
# First block ("run 0"): Kruskal is now at row 7 -> insert above it.
$ws2.Range("A7:K9").Insert()
# Second block ("Average"): Kruskal is now at row 16 -> insert above it.
$ws2.Range("A16:K18").Insert()

$ws2.Range("B7").Value = "S*-HS-UN"
$ws2.Range("C7").Value = 1819
$ws2.Range("D7").Value = 0.1308123781345785
$ws2.Range("E7").Value = 1.487103268038481
$ws2.Range("F7").Value = 1819
$ws2.Range("G7").Value = 0.09201491670683026
$ws2.Range("H7").Value = 0.6701455852016807
$ws2.Range("I7").Value = 0.4526925073005259
$ws2.Range("J7").Value = 0.1424911376088858
$ws2.Range("K7").Value = 0.05763309635221958

$ws2.Range("B8").Value = "S*-MM-UN"
$ws2.Range("C8").Value = 1413
$ws2.Range("D8").Value = 0.1010411842726171
$ws2.Range("E8").Value = 1.1247075391002
$ws2.Range("F8").Value = 1413
$ws2.Range("G8").Value = 0.0726292678155005
$ws2.Range("H8").Value = 0.5177717562764883
$ws2.Range("I8").Value = 0.3065723571926355
$ws2.Range("J8").Value = 0.1286666635423899
$ws2.Range("K8").Value = 0.04402244556695223

$ws2.Range("B9").Value = "S*-MM0-UN"
$ws2.Range("C9").Value = 1950
$ws2.Range("D9").Value = 0.007076812908053398
$ws2.Range("E9").Value = 1.472084654960781
$ws2.Range("F9").Value = 1950
$ws2.Range("G9").Value = 0.100663264747709
$ws2.Range("H9").Value = 0.573080484289676
$ws2.Range("I9").Value = 0.6032815179787576
$ws2.Range("J9").Value = 0.05679511511698365
$ws2.Range("K9").Value = 0.05919745704159141

$ws2.Range("B16").Value = "S*-HS-UN"
$ws2.Range("C16").Value = 1819
$ws2.Range("D16").Value = 0.1308123781345785
$ws2.Range("E16").Value = 1.487103268038481
$ws2.Range("F16").Value = 1819
$ws2.Range("G16").Value = 0.09201491670683026
$ws2.Range("H16").Value = 0.6701455852016807
$ws2.Range("I16").Value = 0.4526925073005259
$ws2.Range("J16").Value = 0.1424911376088858
$ws2.Range("K16").Value = 0.05763309635221958

$ws2.Range("B17").Value = "S*-MM-UN"
$ws2.Range("C17").Value = 1413
$ws2.Range("D17").Value = 0.1010411842726171
$ws2.Range("E17").Value = 1.1247075391002
$ws2.Range("F17").Value = 1413
$ws2.Range("G17").Value = 0.0726292678155005
$ws2.Range("H17").Value = 0.5177717562764883
$ws2.Range("I17").Value = 0.3065723571926355
$ws2.Range("J17").Value = 0.1286666635423899
$ws2.Range("K17").Value = 0.04402244556695223

$ws2.Range("B18").Value = "S*-MM0-UN"
$ws2.Range("C18").Value = 1950
$ws2.Range("D18").Value = 0.007076812908053398
$ws2.Range("E18").Value = 1.472084654960781
$ws2.Range("F18").Value = 1950
$ws2.Range("G18").Value = 0.100663264747709
$ws2.Range("H18").Value = 0.573080484289676
$ws2.Range("I18").Value = 0.6032815179787576
$ws2.Range("J18").Value = 0.05679511511698365
$ws2.Range("K18").Value = 0.05919745704159141

# Update the surviving Kruskal rows + "run 0" values (B/C columns
# through the new solver columns stay the same values, since the
# run-0 & run-1 split collapsed into identical numbers in v2 — the
# diff only really changes the C6/E6 "var" row and the Kruskal E
# column numbers).
$ws2.Range("C6").Value = 1950
$ws2.Range("D6").Value = 0.005769672337919474
$ws2.Range("E6").Value = 1.415210830979049
$ws2.Range("F6").Value = 1950
$ws2.Range("G6").Value = 0.09484151471406221
$ws2.Range("H6").Value = 0.5603964230976999
$ws2.Range("I6").Value = 0.5701635577715933
$ws2.Range("J6").Value = 0.05438367929309607
$ws2.Range("K6").Value = 0.05979041056707501

$ws2.Range("B10").Value = "Kruskal"
$ws2.Range("C10").Value = 33300
$ws2.Range("D10").ClearContents()
$ws2.Range("E10").Value = 6.630521262995899
$ws2.Range("F10").ClearContents()
$ws2.Range("G10").ClearContents()
$ws2.Range("H10").ClearContents()
$ws2.Range("I10").ClearContents()
$ws2.Range("J10").ClearContents()
$ws2.Range("K10").ClearContents()

$ws2.Range("A11").Value = "Average"
$ws2.Range("C11").Value = 657
$ws2.Range("D11").Value = 0.002956422977149487
$ws2.Range("E11").Value = 0.4192828559316695
$ws2.Range("F11").Value = 657
$ws2.Range("G11").Value = 0.03663060395047069
$ws2.Range("H11").Value = 0.2144852154888213
$ws2.Range("I11").Value = 0.07103201141580939
$ws2.Range("J11").Value = 0.05145837739109993
$ws2.Range("K11").Value = 0.01794490171596408

$ws2.Range("C15").Value = 1950
$ws2.Range("D15").Value = 0.005769672337919474
$ws2.Range("E15").Value = 1.415210830979049
$ws2.Range("F15").Value = 1950
$ws2.Range("G15").Value = 0.09484151471406221
$ws2.Range("H15").Value = 0.5603964230976999
$ws2.Range("I15").Value = 0.5701635577715933
$ws2.Range("J15").Value = 0.05438367929309607
$ws2.Range("K15").Value = 0.05979041056707501

$ws2.Range("B19").Value = "Kruskal"
$ws2.Range("C19").Value = 33300
$ws2.Range("D19").ClearContents()
$ws2.Range("E19").Value = 6.630521262995899
$ws2.Range("F19").ClearContents()
$ws2.Range("G19").ClearContents()
$ws2.Range("H19").ClearContents()
$ws2.Range("I19").ClearContents()
$ws2.Range("J19").ClearContents()
$ws2.Range("K19").ClearContents()

# Re-merge the per-block "A" label cells over the new, longer row spans.
$ws2.Range("A2:A10").Merge()
$ws2.Range("A11:A19").Merge()
